# Fruta / hortaliza, semanal
# Insert a new weekly observation row at row 279 (pushing the existing
# rows 279-289 down to 280-290) for "Vega Central Mapocho de Santiago - Mango".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 279:289 down to 280:290, leaving a fresh blank row at 279.
$ws.Rows.Item(279).Insert()

# Populate the new row 279 with the new weekly record.
$ws.Cells.Item(279, 1).Value = 9
$ws.Cells.Item(279, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(279, 3).Value = "Metropolitana"
$ws.Cells.Item(279, 4).Value = 44516
$ws.Cells.Item(279, 5).Value = 13
$ws.Cells.Item(279, 6).Value = "Fruta"
$ws.Cells.Item(279, 7).Value = 100108
$ws.Cells.Item(279, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(279, 9).Value = 100108002
$ws.Cells.Item(279, 10).Value = "Mango"
$ws.Cells.Item(279, 11).Value = "Sin especificar"
$ws.Cells.Item(279, 12).Value = "Primera"
$ws.Cells.Item(279, 13).Value = 660
$ws.Cells.Item(279, 14).Value = 6000
$ws.Cells.Item(279, 15).Value = 6500
$ws.Cells.Item(279, 16).Value = 6265
$ws.Cells.Item(279, 17).Value = "`$/bandeja 4 kilos"
$ws.Cells.Item(279, 18).Value = "Brasil"
$ws.Cells.Item(279, 19).Value = 1566
$ws.Cells.Item(279, 20).Value = 4
